$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 61441.59
$ws.Range("I64").Value = 102420.1
$ws.Range("J64").Value = 2900.8572
$ws.Range("K64").Value = 102420.1
$ws.Range("L64").Value = 2900.8572
$ws.Range("M64").Value = -102172.1
$ws.Range("N64").Value = -3396.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 61441.59
$ws.Range("I67").Value = 102420.1
$ws.Range("J67").Value = 2900.8572
$ws.Range("K67").Value = 102420.1
$ws.Range("L67").Value = 2900.8572
$ws.Range("M67").Value = -101562.1
$ws.Range("N67").Value = -4616.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4573.625
$ws.Range("I106").Value = 4764.8335
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 4764.8335
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -4133.8335
$ws.Range("N106").Value = -5262

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2432.7778
$ws.Range("I116").Value = 1873.5714
$ws.Range("J116").Value = 4390
$ws.Range("K116").Value = 1873.5714
$ws.Range("L116").Value = 4390
$ws.Range("M116").Value = 1568.4286
$ws.Range("N116").Value = -11274

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2947
$ws.Range("J129").Value = 1375.1936
$ws.Range("L129").Value = 4125.5808
$ws.Range("N129").Value = -14125.5808

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 974.36365
$ws.Range("I137").Value = 968.7193
$ws.Range("J137").Value = 1010.1111
$ws.Range("K137").Value = 2906.1579
$ws.Range("L137").Value = 3030.3333
$ws.Range("M137").Value = -356.1579000000002
$ws.Range("N137").Value = -8130.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1519.5333
$ws.Range("I138").Value = 983.907
$ws.Range("J138").Value = 2874.353
$ws.Range("K138").Value = 2951.721
$ws.Range("L138").Value = 8623.059000000001
$ws.Range("M138").Value = 2188.279
$ws.Range("N138").Value = -18903.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4922.22
$ws.Range("I32").Value = 3774.9333
$ws.Range("J32").Value = 15247.8
$ws.Range("K32").Value = 3774.9333
$ws.Range("L32").Value = 15247.8
$ws.Range("M32").Value = -3487.9333
$ws.Range("N32").Value = -15821.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 432.6579
$ws.Range("I74").Value = 416.27274
$ws.Range("J74").Value = 540.8
$ws.Range("K74").Value = 416.27274
$ws.Range("L74").Value = 540.8
$ws.Range("M74").Value = 457.72726
$ws.Range("N74").Value = -2288.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 432.6579
$ws.Range("I77").Value = 416.27274
$ws.Range("J77").Value = 540.8
$ws.Range("K77").Value = 2081.3637
$ws.Range("L77").Value = 2704
$ws.Range("M77").Value = 2286.6363
$ws.Range("N77").Value = -11440

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4290.0527
$ws.Range("I132").Value = 4285.1875
$ws.Range("J132").Value = 4316
$ws.Range("K132").Value = 12855.5625
$ws.Range("L132").Value = 12948
$ws.Range("M132").Value = -10325.5625
$ws.Range("N132").Value = -18008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 39630.5
$ws.Range("J133").Value = 39630.5
$ws.Range("L133").Value = 39630.5
$ws.Range("N133").Value = -44690.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1090.1666
$ws.Range("I16").Value = 758.7
$ws.Range("K16").Value = 758.7
$ws.Range("M16").Value = -471.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24634.396
$ws.Range("I31").Value = 1463.8684
$ws.Range("K31").Value = 1463.8684
$ws.Range("M31").Value = -1168.8684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 24634.396
$ws.Range("I34").Value = 1463.8684
$ws.Range("K34").Value = 1463.8684
$ws.Range("M34").Value = -1261.8684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1100.3334
$ws.Range("I58").Value = 917.3333
$ws.Range("K58").Value = 917.3333
$ws.Range("M58").Value = -714.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1090.1666
$ws.Range("I113").Value = 758.7
$ws.Range("K113").Value = 758.7
$ws.Range("M113").Value = 1411.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2986.878
$ws.Range("I132").Value = 2550.276
$ws.Range("J132").Value = 4042
$ws.Range("K132").Value = 7650.828
$ws.Range("L132").Value = 12126
$ws.Range("M132").Value = -5120.828
$ws.Range("N132").Value = -17186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1100.3334
$ws.Range("I136").Value = 917.3333
$ws.Range("K136").Value = 2751.9999
$ws.Range("M136").Value = -201.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 650.17645
$ws.Range("I23").Value = 545.75
$ws.Range("J23").Value = 682.3077
$ws.Range("K23").Value = 1637.25
$ws.Range("L23").Value = 2046.9231
$ws.Range("M23").Value = -1402.25
$ws.Range("N23").Value = -2516.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 83714.164
$ws.Range("J98").Value = 143127.14
$ws.Range("L98").Value = 429381.42
$ws.Range("N98").Value = -432377.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2721.111
$ws.Range("I132").Value = 766.6667
$ws.Range("K132").Value = 6900.0003
$ws.Range("M132").Value = -4370.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 63180.293
$ws.Range("I70").Value = 87430.664
$ws.Range("J70").Value = 4979.4
$ws.Range("K70").Value = 87430.664
$ws.Range("L70").Value = 4979.4
$ws.Range("M70").Value = -87160.664
$ws.Range("N70").Value = -5519.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 63180.293
$ws.Range("I73").Value = 87430.664
$ws.Range("J73").Value = 4979.4
$ws.Range("K73").Value = 87430.664
$ws.Range("L73").Value = 4979.4
$ws.Range("M73").Value = -86494.664
$ws.Range("N73").Value = -6851.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2395.12
$ws.Range("I102").Value = 1884
$ws.Range("K102").Value = 1884
$ws.Range("M102").Value = -262

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4524.364
$ws.Range("I132").Value = 4536.706
$ws.Range("J132").Value = 4482.4
$ws.Range("K132").Value = 13610.118
$ws.Range("L132").Value = 13447.2
$ws.Range("M132").Value = -11080.118
$ws.Range("N132").Value = -18507.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 325667.1
$ws.Range("I55").Value = 1136652.4
$ws.Range("J55").Value = 1272.96
$ws.Range("K55").Value = 1136652.4
$ws.Range("L55").Value = 1272.96
$ws.Range("M55").Value = -1136479.4
$ws.Range("N55").Value = -1618.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2842.476
$ws.Range("I132").Value = 2718.5405
$ws.Range("J132").Value = 3759.6
$ws.Range("K132").Value = 8155.6215
$ws.Range("L132").Value = 11278.8
$ws.Range("M132").Value = -5625.6215
$ws.Range("N132").Value = -16338.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 923.1587500000001
$ws.Range("I136").Value = 760.1818
$ws.Range("J136").Value = 2043.625
$ws.Range("K136").Value = 2280.5454
$ws.Range("L136").Value = 6130.875
$ws.Range("M136").Value = 269.4546
$ws.Range("N136").Value = -11230.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 47952.76
$ws.Range("I107").Value = 216.86667
$ws.Range("J107").Value = 167292.5
$ws.Range("K107").Value = 650.60001
$ws.Range("L107").Value = 501877.5
$ws.Range("M107").Value = 1269.39999
$ws.Range("N107").Value = -505717.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1489.05
$ws.Range("I126").Value = 1447.3572
$ws.Range("J126").Value = 1586.3334
$ws.Range("K126").Value = 4342.071599999999
$ws.Range("L126").Value = 4759.0002
$ws.Range("M126").Value = -1872.071599999999
$ws.Range("N126").Value = -9699.0002
